# Adressed more rmTypes in pathExport. Few missing
#
# Updates the FLAT_Paths sheet (column A, the list of FLAT paths used for the
# data-validation dropdowns on the "Mapping CSV2openEHR" sheet) so several
# rmType paths get the correct "|value" / "|code" suffixes, matching the
# order used elsewhere in the mapping (value before code).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FLAT_Paths")

$updates = @{
    2  = "natars_tzusatz/category|value"
    3  = "natars_tzusatz/category|code"
    7  = "natars_tzusatz/context/setting|value"
    8  = "natars_tzusatz/context/setting|code"
    10 = "natars_tzusatz/prevention_of_falls_network_europe_profane/sturz|value"
    11 = "natars_tzusatz/prevention_of_falls_network_europe_profane/sturz|code"
    13 = "natars_tzusatz/prevention_of_falls_network_europe_profane/häufigkeit_des_sturzes|value"
    15 = "natars_tzusatz/prevention_of_falls_network_europe_profane/verletzung_erlitten|value"
    16 = "natars_tzusatz/prevention_of_falls_network_europe_profane/verletzung_erlitten|code"
    18 = "natars_tzusatz/prevention_of_falls_network_europe_profane/bruch|value"
    19 = "natars_tzusatz/prevention_of_falls_network_europe_profane/bruch|code"
    45 = "natars_tzusatz/schmerzerfassung/schmerzstärke/schmerzstärke_-_vas/beliebiges_ereignis<<index>>/vas-wert|value"
    56 = "natars_tzusatz/schmerzerfassung/bewegungseinschränkung_durch_die_schmerzen/bewegungseinschränkung/beliebiges_ereignis<<index>>/vas-wert|value"
    71 = "natars_tzusatz/body_mass_index/methode|value"
    72 = "natars_tzusatz/body_mass_index/methode|code"
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
